$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 665
$ws1.Cells.Item(7, 6).Value = 1252
$ws1.Cells.Item(9, 6).Value = 2395
$ws1.Cells.Item(10, 6).Value = 841
$ws1.Cells.Item(11, 6).Value = 18300
$ws1.Cells.Item(12, 6).Value = 41
$ws1.Cells.Item(13, 6).Value = 1758
$ws1.Cells.Item(14, 6).Value = 632
$ws1.Cells.Item(15, 6).Value = 588
$ws1.Cells.Item(16, 6).Value = 301
$ws1.Cells.Item(17, 6).Value = 583
$ws1.Cells.Item(18, 6).Value = 183
$ws1.Cells.Item(19, 6).Value = 181
$ws1.Cells.Item(21, 6).Value = 308
$ws1.Cells.Item(22, 6).Value = 158
$ws1.Cells.Item(23, 6).Value = 78
$ws1.Cells.Item(24, 6).Value = 11

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(7, 6).Value = 105
$ws2.Cells.Item(8, 6).Value = 104
$ws2.Cells.Item(14, 6).Value = 64

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5829
$ws3.Cells.Item(3, 6).Value = 535

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 5829
$ws4.Cells.Item(4, 6).Value = 535
$ws4.Cells.Item(7, 6).Value = 665
$ws4.Cells.Item(12, 6).Value = 1252
$ws4.Cells.Item(17, 6).Value = 2395
$ws4.Cells.Item(18, 6).Value = 841
$ws4.Cells.Item(19, 6).Value = 18301
$ws4.Cells.Item(20, 6).Value = 41
$ws4.Cells.Item(21, 6).Value = 105
$ws4.Cells.Item(22, 6).Value = 104
$ws4.Cells.Item(23, 6).Value = 104
$ws4.Cells.Item(24, 6).Value = 1758
$ws4.Cells.Item(25, 6).Value = 632
$ws4.Cells.Item(27, 6).Value = 588
$ws4.Cells.Item(28, 6).Value = 301
$ws4.Cells.Item(29, 6).Value = 583
$ws4.Cells.Item(30, 6).Value = 183
$ws4.Cells.Item(31, 6).Value = 181
$ws4.Cells.Item(36, 6).Value = 308
$ws4.Cells.Item(38, 6).Value = 64
$ws4.Cells.Item(39, 6).Value = 158
$ws4.Cells.Item(41, 6).Value = 78
$ws4.Cells.Item(45, 6).Value = 11
